$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()
$ws.Range("H132").Value = 1284.3334
$ws.Range("I132").Value = 1037.6818
$ws.Range("K132").Value = 3113.0454
$ws.Range("M132").Value = -583.0454
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 341.8
$ws.Range("I2").Value = 349.75
$ws.Range("J2").Value = 310
$ws.Range("K2").Value = 349.75
$ws.Range("L2").Value = 310
$ws.Range("M2").Value = -236.75
$ws.Range("N2").Value = -536
$ws.Range("H32").Value = 4017.6897
$ws.Range("I32").Value = 4089.75
$ws.Range("K32").Value = 4089.75
$ws.Range("M32").Value = -3802.75
$ws.Range("H74").Value = 2012
$ws.Range("I74").Value = 2012
$ws.Range("K74").Value = 2012
$ws.Range("M74").Value = -1138
$ws.Range("H77").Value = 2012
$ws.Range("I77").Value = 2012
$ws.Range("K77").Value = 10060
$ws.Range("M77").Value = -5692
$ws.Range("H97").Value = 5387.7
$ws.Range("I97").Value = 4782.2856
$ws.Range("J97").Value = 6800.3335
$ws.Range("K97").Value = 4782.2856
$ws.Range("L97").Value = 6800.3335
$ws.Range("M97").Value = -4286.2856
$ws.Range("N97").Value = -7792.3335
$ws.Range("H102").Value = 1853
$ws.Range("I102").Value = 1853
$ws.Range("K102").Value = 1853
$ws.Range("M102").Value = -231
$ws.Range("H116").Value = 341.8
$ws.Range("I116").Value = 349.75
$ws.Range("J116").Value = 310
$ws.Range("K116").Value = 349.75
$ws.Range("L116").Value = 310
$ws.Range("M116").Value = 1944.25
$ws.Range("N116").Value = -4898
$ws.Range("H122").Value = 11059.714
$ws.Range("I122").Value = 11059.714
$ws.Range("K122").Value = 33179.142
$ws.Range("M122").Value = -30729.142
$ws.Range("H134").Value = 28000
$ws.Range("J134").Value = 28000
$ws.Range("L134").Value = 28000
$ws.Range("N134").Value = -38140

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 341.8
$ws.Range("I3").Value = 349.75
$ws.Range("J3").Value = 310
$ws.Range("K3").Value = 349.75
$ws.Range("L3").Value = 310
$ws.Range("M3").Value = -235.75
$ws.Range("N3").Value = -538

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3152.394
$ws.Range("I31").Value = 2462.6924
$ws.Range("J31").Value = 5714.143
$ws.Range("K31").Value = 2462.6924
$ws.Range("L31").Value = 5714.143
$ws.Range("M31").Value = -2167.6924
$ws.Range("N31").Value = -6304.143
$ws.Range("H34").Value = 3152.394
$ws.Range("I34").Value = 2462.6924
$ws.Range("J34").Value = 5714.143
$ws.Range("K34").Value = 2462.6924
$ws.Range("L34").Value = 5714.143
$ws.Range("M34").Value = -2260.6924
$ws.Range("N34").Value = -6118.143
$ws.Range("H86").Value = 7241.636
$ws.Range("I86").Value = 4746.25
$ws.Range("K86").Value = 4746.25
$ws.Range("M86").Value = -3623.25
$ws.Range("H89").Value = 7241.636
$ws.Range("I89").Value = 4746.25
$ws.Range("K89").Value = 23731.25
$ws.Range("M89").Value = -18115.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 1144
$ws.Range("I18").Value = 1152.8
$ws.Range("K18").Value = 3458.4
$ws.Range("M18").Value = -3289.4
$ws.Range("H34").Value = 3751.125
$ws.Range("I34").Value = 3009.75
$ws.Range("J34").Value = 3998.25
$ws.Range("K34").Value = 9029.25
$ws.Range("L34").Value = 11994.75
$ws.Range("M34").Value = -8945.25
$ws.Range("N34").Value = -12162.75
$ws.Range("H139").Value = 2449.6667
$ws.Range("I139").Value = 2449.6667
$ws.Range("K139").Value = 7349.000100000001
$ws.Range("M139").Value = -2209.000100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2744
$ws.Range("J102").Value = 3649.25
$ws.Range("L102").Value = 3649.25
$ws.Range("N102").Value = -6893.25
$ws.Range("H113").Value = 2546.4
$ws.Range("I113").Value = 2433
$ws.Range("K113").Value = 2433
$ws.Range("M113").Value = -263
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H132").Value = 2130.5
$ws.Range("I132").Value = 1822.7333
$ws.Range("J132").Value = 3669.3333
$ws.Range("K132").Value = 5468.199900000001
$ws.Range("L132").Value = 11007.9999
$ws.Range("M132").Value = -2938.199900000001
$ws.Range("N132").Value = -16067.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2775
$ws.Range("I100").Value = 2775
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 2775
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -2234
$ws.Range("N100").ClearContents()
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H136").Value = 1917
$ws.Range("I136").Value = 1917
$ws.Range("K136").Value = 5751
$ws.Range("M136").Value = -3201

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 614.75
$ws.Range("I122").Value = 536.75
$ws.Range("K122").Value = 1610.25
$ws.Range("M122").Value = 839.75
$ws.Range("H132").Value = 3698.138
$ws.Range("I132").Value = 2319.7778
$ws.Range("K132").Value = 6959.3334
$ws.Range("M132").Value = -4429.3334
$ws.Range("H136").Value = 2077.9092
$ws.Range("I136").Value = 1907.125
$ws.Range("J136").Value = 2533.3333
$ws.Range("K136").Value = 5721.375
$ws.Range("L136").Value = 7599.999899999999
$ws.Range("M136").Value = -3171.375
$ws.Range("N136").Value = -12699.9999
